$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ C=35; D=34; E=34; F=37; G=8;  H=7; I=10; J=5; K=10; L=9; M=11; N=7 }
    3 = @{ C=29; D=34; E=27; F=38; G=6;  H=6; I=9;  J=4; K=6;  L=8; M=7;  N=7 }
    4 = @{ C=3;  D=8;  E=6;  F=6;  G=5;  H=6; I=6;  J=2; K=5;  L=8; M=7;  N=5 }
    5 = @{ C=5;  D=2;  E=4;  F=3;  G=6;  H=5; I=6;  J=4; K=7;  L=7; M=8;  N=5 }
    6 = @{ C=3;  D=7;  E=6;  F=6;  G=5;  H=4; I=6;  J=2; K=5;  L=7; M=7;  N=5 }
    7 = @{ C=6;  D=3;  E=5;  F=4;  G=5;  H=4; I=6;  J=4; K=7;  L=7; M=8;  N=5 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
